$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period (fiscal-year) headers shift left by one year, new year appended in H ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates shift left by one year, new date appended in H ---
$ws.Range("D9").Value = "1399-03-07 (8)"
$ws.Range("E9").Value = "1400-02-28 (9)"
$ws.Range("F9").Value = "1401-03-11 (8)"
$ws.Range("G9").Value = "1402-02-30 (7)"
$ws.Range("H9").Value = "1402-02-30"

# --- Row 11: Sales (فروش) ---
$ws.Range("D11").Value = 2693277
$ws.Range("E11").Value = 4031266
$ws.Range("F11").Value = 8468536
$ws.Range("G11").Value = 13691281
$ws.Range("H11").Value = 24645941

# --- Row 12: Cost of goods sold ---
$ws.Range("D12").Value = -2689066
$ws.Range("E12").Value = -3595540
$ws.Range("F12").Value = -7313605
$ws.Range("G12").Value = -11319270
$ws.Range("H12").Value = -19785175

# --- Row 13: Gross profit (loss) ---
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 435726
$ws.Range("F13").Value = 1154931
$ws.Range("G13").Value = 2372011
$ws.Range("H13").Value = 4860766

# --- Row 14: General, administrative and organizational expenses ---
$ws.Range("D14").Value = -251216
$ws.Range("E14").Value = -354183
$ws.Range("F14").Value = -791950
$ws.Range("G14").Value = -1093286
$ws.Range("H14").Value = -1359815

# --- Row 15: Impairment of receivables expense (D15 switches from "-" text to numeric 0) ---
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: Net other operating income (expenses) ---
$ws.Range("D16").Value = 409345
$ws.Range("E16").Value = 241410
$ws.Range("F16").Value = 411446
$ws.Range("G16").Value = 4537
$ws.Range("H16").Value = 883120

# --- Row 17: Operating profit (loss) ---
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 322953
$ws.Range("F17").Value = 774427
$ws.Range("G17").Value = 1283262
$ws.Range("H17").Value = 4384071

# --- Row 18: Financial expenses ---
$ws.Range("D18").Value = -4698
$ws.Range("E18").Value = -4400
$ws.Range("F18").Value = -17293
$ws.Range("G18").Value = -48202
$ws.Range("H18").Value = -15641

# --- Row 19: Net other non-operating income and expenses ---
$ws.Range("D19").Value = 12037
$ws.Range("E19").Value = 23537
$ws.Range("F19").Value = 148974
$ws.Range("G19").Value = 215534
$ws.Range("H19").Value = 141080

# --- Row 20: Net profit (loss) from continuing operations before tax ---
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 342090
$ws.Range("F20").Value = 906108
$ws.Range("G20").Value = 1450594
$ws.Range("H20").Value = 4509510

# --- Row 21: Tax ---
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = -5083
$ws.Range("G21").Value = -92886
$ws.Range("H21").Value = -510870

# --- Row 22: Net profit (loss) from continuing operations ---
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 342090
$ws.Range("F22").Value = 901025
$ws.Range("G22").Value = 1357708
$ws.Range("H22").Value = 3998640

# --- Row 23: Profit (loss) from discontinued operations after tax effect ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: Net profit (loss) ---
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 342090
$ws.Range("F24").Value = 901025
$ws.Range("G24").Value = 1357708
$ws.Range("H24").Value = 3998640

# --- Row 25: Earnings per share after tax ---
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 608
$ws.Range("F25").Value = 493
$ws.Range("G25").Value = 744
$ws.Range("H25").Value = 2190

# --- Row 26: Capital ---
$ws.Range("D26").Value = 213000
$ws.Range("E26").Value = 563000
$ws.Range("F26").Value = 1826000
$ws.Range("G26").Value = 1826000
$ws.Range("H26").Value = 1826000

# --- Row 27: Earnings per share based on latest capital ---
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 187
$ws.Range("F27").Value = 493
$ws.Range("G27").Value = 744
$ws.Range("H27").Value = 2190

Write-Host "edit complete"
